$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(2, 2).Value = "Nemacheilidae (n=18)"
$ws.Cells.Item(2, 3).Value = 8.680848393574308
$ws.Cells.Item(2, 4).Value = 0.003215712573055775
$ws.Cells.Item(2, 5).Value = 0.01607856286527887

$ws.Cells.Item(3, 1).Value = "Nemacheilidae (n=18)"
$ws.Cells.Item(3, 2).Value = "Tilapiinae (n=12)"
$ws.Cells.Item(3, 3).Value = 8.781362007168454
$ws.Cells.Item(3, 4).Value = 0.003043238861600076
$ws.Cells.Item(3, 5).Value = 0.01607856286527887

$ws.Cells.Item(8, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(8, 2).Value = "Haplochrominae (n=3)"
$ws.Cells.Item(8, 3).Value = 1.191176470588232
$ws.Cells.Item(8, 4).Value = 0.2750923766358885
$ws.Cells.Item(8, 5).Value = 0.3563117592532612

$ws.Cells.Item(10, 1).Value = "Haplochrominae (n=3)"
$ws.Cells.Item(10, 2).Value = "Mugilidae (n=1)"
$ws.Cells.Item(10, 3).Value = 0.1999999999999993
$ws.Cells.Item(10, 4).Value = 0.6547208460185774
$ws.Cells.Item(10, 5).Value = 0.7274676066873081

$ws.Cells.Item(11, 1).Value = "Haplochrominae (n=3)"
$ws.Cells.Item(11, 2).Value = "Tilapiinae (n=12)"
$ws.Cells.Item(11, 3).Value = 0.02083333333334281
$ws.Cells.Item(11, 4).Value = 0.8852339144731757
$ws.Cells.Item(11, 5).Value = 0.8852339144731757

$ws.Cells.Item(12, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(12, 2).Value = "Haplochrominae (n=8)"
$ws.Cells.Item(12, 3).Value = 8.425941780821915
$ws.Cells.Item(12, 4).Value = 0.003699050019305536
$ws.Cells.Item(12, 5).Value = 0.03939428548774692

$ws.Cells.Item(13, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(13, 2).Value = "Tilapiinae (n=22)"
$ws.Cells.Item(13, 3).Value = 7.790360501567419
$ws.Cells.Item(13, 4).Value = 0.005252571398366256
$ws.Cells.Item(13, 5).Value = 0.03939428548774692

$ws.Cells.Item(14, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(14, 2).Value = "Mugilidae (n=3)"
$ws.Cells.Item(14, 3).Value = 4.897977941176464
$ws.Cells.Item(14, 4).Value = 0.02688816202716629
$ws.Cells.Item(14, 5).Value = 0.1344408101358315

$ws.Cells.Item(16, 1).Value = "Haplochrominae (n=8)"
$ws.Cells.Item(16, 2).Value = "Poeciliidae (n=2)"
$ws.Cells.Item(16, 3).Value = 3.340909090909086
$ws.Cells.Item(16, 4).Value = 0.06757726305587061
$ws.Cells.Item(16, 5).Value = 0.1784239642790366

$ws.Cells.Item(17, 1).Value = "Mugilidae (n=3)"
$ws.Cells.Item(17, 2).Value = "Nemacheilidae (n=21)"
$ws.Cells.Item(17, 3).Value = 3.201904761904757
$ws.Cells.Item(17, 4).Value = 0.07355255978200786
$ws.Cells.Item(17, 5).Value = 0.1784239642790366

$ws.Cells.Item(19, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(19, 2).Value = "Nemacheilidae (n=21)"
$ws.Cells.Item(19, 3).Value = 1.948608803986701
$ws.Cells.Item(19, 4).Value = 0.1627368442509667
$ws.Cells.Item(19, 5).Value = 0.3051315829705624

$ws.Cells.Item(20, 1).Value = "Nemacheilidae (n=21)"
$ws.Cells.Item(20, 2).Value = "Tilapiinae (n=22)"
$ws.Cells.Item(20, 3).Value = 1.721369539551347
$ws.Cells.Item(20, 4).Value = 0.1895168505762421
$ws.Cells.Item(20, 5).Value = 0.3158614176270702

$ws.Cells.Item(21, 1).Value = "Poeciliidae (n=2)"
$ws.Cells.Item(21, 2).Value = "Tilapiinae (n=22)"
$ws.Cells.Item(21, 3).Value = 1.320000000000007
$ws.Cells.Item(21, 4).Value = 0.2505920506856796
$ws.Cells.Item(21, 5).Value = 0.3758880760285194

$ws.Cells.Item(22, 1).Value = "Mugilidae (n=3)"
$ws.Cells.Item(22, 2).Value = "Tilapiinae (n=22)"
$ws.Cells.Item(22, 3).Value = 1.006993006993
$ws.Cells.Item(22, 4).Value = 0.3156243007353199
$ws.Cells.Item(22, 5).Value = 0.4303967737299818

$ws.Cells.Item(23, 1).Value = "Haplochrominae (n=8)"
$ws.Cells.Item(23, 2).Value = "Mugilidae (n=3)"
$ws.Cells.Item(23, 3).Value = 0.6666666666666714
$ws.Cells.Item(23, 4).Value = 0.4142161782425236
$ws.Cells.Item(23, 5).Value = 0.4779417441259887

$ws.Cells.Item(24, 1).Value = "Nemacheilidae (n=21)"
$ws.Cells.Item(24, 2).Value = "Poeciliidae (n=2)"
$ws.Cells.Item(24, 3).Value = 0.7619047619047592
$ws.Cells.Item(24, 4).Value = 0.3827330888852269
$ws.Cells.Item(24, 5).Value = 0.4779417441259887

$ws.Cells.Item(25, 1).Value = "Haplochrominae (n=8)"
$ws.Cells.Item(25, 2).Value = "Tilapiinae (n=22)"
$ws.Cells.Item(25, 3).Value = 0.3717008797653989
$ws.Cells.Item(25, 4).Value = 0.542078599605877
$ws.Cells.Item(25, 5).Value = 0.5807984995777253

$ws.Cells.Item(26, 1).Value = "Cyprinidae (n=64)"
$ws.Cells.Item(26, 2).Value = "Poeciliidae (n=2)"
$ws.Cells.Item(26, 3).Value = 0.05037313432833912
$ws.Cells.Item(26, 4).Value = 0.822415259699475
$ws.Cells.Item(26, 5).Value = 0.822415259699475
